$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("种类" / item-type) held "传奇道具" ("legendary item") for both
# data rows; the card-number/type field is simplified to just "道具".
$ws.Range("D2").Value = "道具"
$ws.Range("D3").Value = "道具"

$ws.Range("I14").Select() | Out-Null
